$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows at position 6, pushing existing rows 6..105 down to 8..107
$ws.Rows("6:7").Insert()

# New row 6 (Fecha 45083, Especial, Volumen 40, Precios 20000/20000/20000, Precio/Kg 1111)
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C6").Value = "Los Lagos"
$ws.Range("D6").Value = 45083
$ws.Range("E6").Value = 10
$ws.Range("F6").Value = 100112043
$ws.Range("G6").Value = "Pepino dulce"
$ws.Range("H6").Value = "Cultivar IV Región"
$ws.Range("I6").Value = "Especial"
$ws.Range("J6").Value = 40
$ws.Range("K6").Value = 20000
$ws.Range("L6").Value = 20000
$ws.Range("M6").Value = 20000
$ws.Range("N6").Value = "$/bandeja 18 kilos"
$ws.Range("O6").Value = "Provincia de Limarí"
$ws.Range("P6").Value = 1111
$ws.Range("Q6").Value = 18
$ws.Range("R6").Value = "Hortaliza"

# New row 7 (Fecha 45083, Primera, Volumen 40, Precios 17000/17000/17000, Precio/Kg 944)
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C7").Value = "Los Lagos"
$ws.Range("D7").Value = 45083
$ws.Range("E7").Value = 10
$ws.Range("F7").Value = 100112043
$ws.Range("G7").Value = "Pepino dulce"
$ws.Range("H7").Value = "Cultivar IV Región"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 40
$ws.Range("K7").Value = 17000
$ws.Range("L7").Value = 17000
$ws.Range("M7").Value = 17000
$ws.Range("N7").Value = "$/bandeja 18 kilos"
$ws.Range("O7").Value = "Provincia de Limarí"
$ws.Range("P7").Value = 944
$ws.Range("Q7").Value = 18
$ws.Range("R7").Value = "Hortaliza"

# Match the date number format used by the rest of column D
$ws.Range("D6:D7").NumberFormat = $ws.Range("D8").NumberFormat
